# Begin work on health summary script:
#  - Swap the "One To One" sheet's AutoFilter on Statistic Category from
#    "Species" to "Health" (this also recomputes which rows are hidden).
#  - Re-point the Overstory Sp 1-5 rows' "Handled by Code Base" /
#    "Dataframe Filter" / "Script Function Name" cells at the new
#    forest_calcs.top5_ov_sp_level helper.
#  - Move the frozen-pane scroll position / active selection to where the
#    author was working (top of sheet, cell F97).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("One To One")

# --- Re-apply the AutoFilter for column A (Statistic Category) ----------
# Switching the single discrete-value filter from "Species" to "Health"
# also updates every row's Hidden state to match the new criteria.
$ws.Range("A1:M236").AutoFilter(1, @("Health"), 7)

# --- Script Function Name (C) rewrites to forest_calcs.top5_ov_sp_level -
$newFn = "forest_calcs.top5_ov_sp_level"

$funcRows = @(149,150,151,152,153,157,158,159,160,161,165,166,167,168,169,173,174,175,176,177,181,182,183,184,185)
foreach ($r in $funcRows) {
    $ws.Range("C$r").Value = $newFn
}

# --- Newly populated Handled by Code Base / Script Function Name /
#     Dataframe Filter cells for the rows that were blank before ----------
$newRows = @{
    154 = "tree_table['TR_SP'] == Overstory Sp 1"
    155 = "tree_table['TR_SP'] == Overstory Sp 1"
    156 = "None"
    162 = "tree_table['TR_SP'] == Overstory Sp 2"
    163 = "tree_table['TR_SP'] == Overstory Sp 2"
    164 = "None"
    170 = "tree_table['TR_SP'] == Overstory Sp 3"
    171 = "tree_table['TR_SP'] == Overstory Sp 3"
    172 = "None"
    178 = "tree_table['TR_SP'] == Overstory Sp 4"
    179 = "tree_table['TR_SP'] == Overstory Sp 4"
    180 = "None"
    186 = "tree_table['TR_SP'] == Overstory Sp 5"
    187 = "tree_table['TR_SP'] == Overstory Sp 5"
}

foreach ($r in $newRows.Keys) {
    $ws.Range("B$r").Value = "Yes"
    $ws.Range("C$r").Value = $newFn
    $ws.Range("D$r").Value = $newRows[$r]
}

# Rows 164/172/180 previously held styled-but-empty B/C/D cells; clear that
# formatting back to Normal now that they carry real values.
foreach ($r in @(164,172,180)) {
    $ws.Range("B$r`:D$r").Style = "Normal"
}

# Row 148 gains a "None" Dataframe Filter value.
$ws.Range("D148").Value = "None"

# --- Restore the author's scroll/selection position ----------------------
$ws.Range("F97").Select()
